# ---------------------------------------------------------------------------
# Edit: insert a new "2022-Q3" sheet (with its summary data) into the
# workbook, right after "总计" and before "2022-Q2".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3 --
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# Copy the number/border/bold formatting used by the other index cells in
# column A (e.g. A6, the former "2020-Q4" row) onto the newly inserted A2.
$summary.Range("A6").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2, 1).Value = 1
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 21
$summary.Cells.Item(2, 4).Value = 1.34

# ---- 2. Insert the brand-new "2022-Q3" detail sheet, right before the
#         existing "2022-Q2" sheet (which is currently the 2nd sheet) -------
$anchor = $wb.Worksheets.Item(2)
$new = $wb.Worksheets.Add($anchor)
$new.Name = "2022-Q3"


# Header row (row 1) labels, columns B..H
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# Detail rows (columns B..H); column A will hold the 0-based row index and
# column H the integer rank, both stored as real numbers like in the other
# sheets. Everything else (fund code / name / size / position values) is
# stored as text, matching the source data format used throughout the
# workbook.
$data = @(
    @("002692", "富国创新科技混合A", "31.88", "86.28", "2.64", "0.8416", 8),
    @("001070", "建信信息产业股票A", "8.92", "91.60", "3.07", "0.2738", 6),
    @("002067", "诺安精选回报灵活配置混合", "1.44", "43.90", "3.13", "0.0451", 7),
    @("512040", "富国中证价值ETF", "2.86", "99.15", "1.24", "0.0355", 5),
    @("000270", "建信灵活配置混合", "2.27", "94.21", "1.05", "0.0238", 1),
    @("002145", "诺安景鑫灵活配置混合", "0.50", "77.19", "3.92", "0.0196", 9),
    @("004194", "招商中证1000指数增强A", "1.56", "92.06", "1.10", "0.0172", 5),
    @("011120", "富国创新科技混合C", "0.64", "86.28", "2.64", "0.0169", 8),
    @("004195", "招商中证1000指数增强C", "1.09", "92.06", "1.10", "0.0120", 5),
    @("013242", "北信瑞丰优势行业股票", "0.69", "92.79", "1.63", "0.0112", 6),
    @("010307", "西藏东财信息产业精选混合A", "0.17", "84.91", "3.74", "0.0064", 8),
    @("165522", "信诚中证TMT产业主题指数（LOF）A", "0.49", "94.08", "1.27", "0.0062", 6),
    @("014246", "大摩现代服务业混合A", "0.17", "66.96", "3.37", "0.0057", 8),
    @("410009", "华富量子生命力混合", "0.11", "93.93", "4.28", "0.0047", 9),
    @("165524", "信诚中证智能家居指数（LOF）A", "0.35", "91.70", "1.22", "0.0043", 7),
    @("010308", "西藏东财信息产业精选混合C", "0.08", "84.91", "3.74", "0.0030", 8),
    @("002952", "建信多因子量化股票", "0.09", "91.26", "3.19", "0.0029", 8),
    @("014247", "大摩现代服务业混合C", "0.06", "66.96", "3.37", "0.0020", 8),
    @("014863", "建信信息产业股票C", "0.06", "91.60", "3.07", "0.0018", 6),
    @("013084", "信诚中证智能家居指数（LOF）C", "0.15", "91.70", "1.22", "0.0018", 7),
    @("013122", "信诚中证TMT产业主题指数（LOF）C", "0.03", "94.08", "1.27", "0.0004", 6)
)

# Force columns B..G to be treated as text before writing, so that values
# such as fund codes ("002692") and percentages ("86.28") keep their
# leading zeros / exact textual representation instead of being coerced to
# numbers by Excel's automatic type detection.
$new.Range("B1:G" + (1 + $data.Count)).NumberFormat = "@"

for ($col = 0; $col -lt $headers.Count; $col++) {
    $new.Cells.Item(1, $col + 2).Value = $headers[$col]
}

$r = 2
foreach ($row in $data) {
    $new.Cells.Item($r, 1).Value = $r - 2
    for ($col = 0; $col -lt 6; $col++) {
        $new.Cells.Item($r, $col + 2).Value = $row[$col]
    }
    $new.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---- 3. Re-apply the workbook's standard header / index-column style ------
# (bold, centered, thin-bordered) to row 1 and column A, matching the style
# already used on every other sheet, and strip the incidental "text number
# format" styling that leaked onto the data cells above.
$summary.Range("B1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A3").Copy()
$new.Range("A2:A" + (1 + $data.Count)).PasteSpecial(-4122)
$new.Range("B2:G" + (1 + $data.Count)).ClearFormats()

$excel.CutCopyMode = $false

# ---- 4. Restore the originally-active sheet ("2020-Q4" was the active tab
#         before this edit) so that adding the new sheet doesn't change the
#         workbook's view/selection state.
$wb.Worksheets.Item("2020-Q4").Activate()
